# Apply the "cryptos list" refresh: updated Price (col D) and
# Volume(1h) (col E) figures for each coin row, plus a reshuffle of
# a handful of rows (40-48) whose Coin/Link/Price/Volume values moved
# to a different row after the ranking changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.060.87'
$ws.Range('E2').Value = '  +0.40%  '
# Row 3
$ws.Range('D3').Value = '3.842.85'
$ws.Range('E3').Value = '  +3.90%  '
# Row 4
$ws.Range('E4').Value = '  +0.19%  '
# Row 5
$ws.Range('D5').Value = '''412.30'
$ws.Range('E5').Value = '  -1.69%  '
# Row 6
$ws.Range('D6').Value = '''131.93'
$ws.Range('E6').Value = '  +1.07%  '
# Row 7
$ws.Range('D7').Value = '3.830.34'
$ws.Range('E7').Value = '  +3.81%  '
# Row 8
$ws.Range('D8').Value = '''0.618'
# Row 9
$ws.Range('E9').Value = '  +0.06%  '
# Row 10
$ws.Range('D10').Value = '''0.739'
$ws.Range('E10').Value = '  -3.54%  '
# Row 11
$ws.Range('D11').Value = '''0.171'
$ws.Range('E11').Value = '  -6.32%  '
# Row 12
$ws.Range('D12').Value = '''0.0000376'
$ws.Range('E12').Value = '  -6.15%  '
# Row 13
$ws.Range('D13').Value = '''41.02'
$ws.Range('E13').Value = '  -4.92%  '
# Row 14
$ws.Range('D14').Value = '4.450.10'
# Row 15
$ws.Range('D15').Value = '''10.04'
$ws.Range('E15').Value = '  -6.23%  '
# Row 16
$ws.Range('D16').Value = '''15.22'
$ws.Range('E16').Value = '  +15.82%  '
# Row 17
$ws.Range('D17').Value = '3.842.68'
$ws.Range('E17').Value = '  +3.74%  '
# Row 18
$ws.Range('E18').Value = '  -1.14%  '
# Row 19
$ws.Range('D19').Value = '''19.59'
$ws.Range('E19').Value = '  -5.07%  '
# Row 20
$ws.Range('D20').Value = '67.453.93'
$ws.Range('E20').Value = '  +0.97%  '
# Row 21
$ws.Range('E21').Value = '  -4.62%  '
# Row 22
$ws.Range('D22').Value = '''414.97'
$ws.Range('E22').Value = '  -6.58%  '
# Row 23
$ws.Range('D23').Value = '''14.61'
$ws.Range('E23').Value = '  -11.21%  '
# Row 24
$ws.Range('D24').Value = '''85.93'
$ws.Range('E24').Value = '  -4.92%  '
# Row 25
$ws.Range('E25').Value = '  -2.43%  '
# Row 26
$ws.Range('D26').Value = '''36.83'
$ws.Range('E26').Value = '  -1.98%  '
# Row 27
$ws.Range('E27').Value = '  +14.28%  '
# Row 28
$ws.Range('D28').Value = '''3.14'
$ws.Range('E28').Value = '  -5.55%  '
# Row 29
$ws.Range('E29').Value = '  -7.17%  '
# Row 30
$ws.Range('D30').Value = '''683.11'
$ws.Range('E30').Value = '  +4.60%  '
# Row 31
$ws.Range('E31').Value = '  -1.84%  '
# Row 32
$ws.Range('E32').Value = '  -3.76%  '
# Row 33
$ws.Range('E33').Value = '  -1.24%  '
# Row 34
$ws.Range('E34').Value = '  -1.08%  '
# Row 35
$ws.Range('E35').Value = '  -8.15%  '
# Row 36
$ws.Range('D36').Value = '''39.08'
$ws.Range('E36').Value = '  -6.02%  '
# Row 37
$ws.Range('D37').Value = '0.0₃0804'
$ws.Range('E37').Value = '  +8.42%  '
# Row 38
$ws.Range('E38').Value = '  +0.03%  '
# Row 39
$ws.Range('D39').Value = '''55.25'
$ws.Range('E39').Value = '  -3.37%  '
# Row 40
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = '''3.09'
$ws.Range('E40').Value = '  +0.31%  '
# Row 41
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0461'
$ws.Range('E41').Value = '  -6.89%  '
# Row 42
$ws.Range('D42').Value = '''0.998'
$ws.Range('E42').Value = '  +0.13%  '
# Row 43
$ws.Range('D43').Value = '''0.137'
$ws.Range('E43').Value = '  -10.07%  '
# Row 44
$ws.Range('D44').Value = '''148.98'
$ws.Range('E44').Value = '  -0.08%  '
# Row 45
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''4.48'
$ws.Range('E45').Value = '  +3.55%  '
# Row 46
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = '''3.34'
$ws.Range('E46').Value = '  -2.73%  '
# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''3.17'
$ws.Range('E47').Value = '  +17.92%  '
# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''26.80'
$ws.Range('E48').Value = '  -10.65%  '
# Row 49
$ws.Range('D49').Value = '''2.09'
$ws.Range('E49').Value = '  -1.33%  '
# Row 50
$ws.Range('D50').Value = '''2.84'
$ws.Range('E50').Value = '  -1.95%  '
# Row 51
$ws.Range('E51').Value = '  -2.19%  '
